$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add I1 ("I0") and J1 ("IF"), matching H1's formatting ---
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# --- Data rows 2-14: I = 1 (constant), J = same value as H (IP) ---
$data = @{
    2  = 2
    3  = 6
    4  = 4
    5  = 6
    6  = 5
    7  = 5
    8  = 4
    9  = 6
    10 = 4
    11 = 3
    12 = 5
    13 = 5
    14 = 3
}

foreach ($r in $data.Keys) {
    $ws.Range("I$r").Value = 1
    $ws.Range("J$r").Value = $data[$r]
}

Write-Host "I0/IF columns added"
